# Generate Report for Archive
$wb = $excel.ActiveWorkbook

# Update status text "Ready for handoff" -> "In Translation" on every sheet
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

# Narrow specific status columns (they were sized to fit "Ready for handoff";
# shrink them to fit the shorter "In Translation" text). Target OOXML column
# width is ~13.41 chars; ColumnWidth = 12.5 is what lands there.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
